$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rearrange header row and data to reflect reading the header row
# to find which columns hold what data.
$ws.Range("A1").Value = "Bunk"
$ws.Range("A2").Value = "B1"

$ws.Range("C1").Value = "Name"
$ws.Range("C2").Value = "Cooper"

$ws.Range("D1").Value = "On Time"
$ws.Range("D2").Value = 2

$ws.Range("B1").Value = "Absent"
$ws.Range("B2").Value = 0

$ws.Range("E1").Value = "Late"
$ws.Range("E2").Value = 3

$ws.Range("B1:B1048576").Select()
